# Add two new columns, I ("I0") and J ("IF"), to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers (copy style from the existing H1 header cell so formatting matches).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2-25 for columns I and J.
$iValues = @{
    2 = 4; 3 = 1; 4 = 1; 5 = 1; 6 = 1; 7 = 1; 8 = 1; 9 = 1; 10 = 1;
    11 = 1; 12 = 1; 13 = 1; 14 = 1; 15 = 1; 16 = 1; 17 = 1; 18 = 1; 19 = 1;
    20 = 1; 21 = 1; 22 = 1; 23 = 1; 24 = 1; 25 = 1
}
$jValues = @{
    2 = 7; 3 = 5; 4 = 4; 5 = 4; 6 = 5; 7 = 4; 8 = 6; 9 = 6; 10 = 6;
    11 = 4; 12 = 5; 13 = 7; 14 = 6; 15 = 6; 16 = 7; 17 = 4; 18 = 5; 19 = 5;
    20 = 7; 21 = 6; 22 = 1; 23 = 5; 24 = 4; 25 = 2
}

for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = $iValues[$r]
    $ws.Cells.Item($r, 10).Value = $jValues[$r]
}
